# Auto-applies numeric value updates to the profit-calculation sheets
# as captured by the upstream data-refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5095
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 5095
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 5095
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -5591
$ws.Range("H67").Value = 5095
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 5095
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 5095
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -6811
$ws.Range("H98").Value = 1285
$ws.Range("I98").Value = 1285
$ws.Range("K98").Value = 1285
$ws.Range("M98").Value = 213
$ws.Range("H116").Value = 4561.968
$ws.Range("J116").Value = 4276.88
$ws.Range("L116").Value = 4276.88
$ws.Range("N116").Value = -11160.88
$ws.Range("H122").Value = 1285
$ws.Range("I122").Value = 1285
$ws.Range("K122").Value = 3855
$ws.Range("M122").Value = -1405
$ws.Range("H132").Value = 5892.4736
$ws.Range("I132").Value = 4527.5586
$ws.Range("K132").Value = 13582.6758
$ws.Range("M132").Value = -11052.6758
$ws.Range("H135").Value = 5180.4375
$ws.Range("I135").Value = 5173.9165
$ws.Range("J135").Value = 5200
$ws.Range("K135").Value = 46565.2485
$ws.Range("L135").Value = 46800
$ws.Range("M135").Value = -44030.2485
$ws.Range("N135").Value = -51870
$ws.Range("H138").Value = 5163.8447
$ws.Range("I138").Value = 4098
$ws.Range("J138").Value = 5441.891
$ws.Range("K138").Value = 12294
$ws.Range("L138").Value = 16325.673
$ws.Range("M138").Value = -7154
$ws.Range("N138").Value = -26605.673

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2852.9
$ws.Range("I61").Value = 2739.9473
$ws.Range("J61").Value = 4999
$ws.Range("K61").Value = 2739.9473
$ws.Range("L61").Value = 4999
$ws.Range("M61").Value = -2527.9473
$ws.Range("N61").Value = -5423
$ws.Range("H74").Value = 161486.58
$ws.Range("I74").Value = 187484.33
$ws.Range("K74").Value = 187484.33
$ws.Range("M74").Value = -186610.33
$ws.Range("H77").Value = 161486.58
$ws.Range("I77").Value = 187484.33
$ws.Range("K77").Value = 937421.6499999999
$ws.Range("M77").Value = -933053.6499999999
$ws.Range("H122").Value = 1228.3914
$ws.Range("I122").Value = 1064.4762
$ws.Range("J122").Value = 2949.5
$ws.Range("K122").Value = 3193.4286
$ws.Range("L122").Value = 8848.5
$ws.Range("M122").Value = -743.4286000000002
$ws.Range("N122").Value = -13748.5
$ws.Range("H132").Value = 28212.922
$ws.Range("I132").Value = 33632.418
$ws.Range("K132").Value = 100897.254
$ws.Range("M132").Value = -98367.25399999999
$ws.Range("H136").Value = 2852.9
$ws.Range("I136").Value = 2739.9473
$ws.Range("J136").Value = 4999
$ws.Range("K136").Value = 8219.841899999999
$ws.Range("L136").Value = 14997
$ws.Range("M136").Value = -5669.841899999999
$ws.Range("N136").Value = -20097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 33900
$ws.Range("J6").Value = 33900
$ws.Range("L6").Value = 33900
$ws.Range("N6").Value = -34126
$ws.Range("H94").Value = 6620.8696
$ws.Range("I94").Value = 7501.6665
$ws.Range("J94").Value = 3450
$ws.Range("K94").Value = 7501.6665
$ws.Range("L94").Value = 3450
$ws.Range("M94").Value = -7050.6665
$ws.Range("N94").Value = -4352
$ws.Range("H107").Value = 2559.6667
$ws.Range("I107").Value = 2491.2856
$ws.Range("J107").Value = 2799
$ws.Range("K107").Value = 2491.2856
$ws.Range("L107").Value = 2799
$ws.Range("M107").Value = -571.2856000000002
$ws.Range("N107").Value = -6639
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H134").Value = 1746.9678
$ws.Range("I134").Value = 1660.5862
$ws.Range("K134").Value = 4981.7586
$ws.Range("M134").Value = -2446.7586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1882.4546
$ws.Range("I22").Value = 280
$ws.Range("J22").Value = 2238.5557
$ws.Range("K22").Value = 280
$ws.Range("L22").Value = 2238.5557
$ws.Range("M22").Value = 70
$ws.Range("N22").Value = -2938.5557
$ws.Range("H47").Value = 38333
$ws.Range("I47").Value = 35000
$ws.Range("J47").Value = 39999.5
$ws.Range("K47").Value = 35000
$ws.Range("L47").Value = 39999.5
$ws.Range("M47").Value = -34434
$ws.Range("N47").Value = -41131.5
$ws.Range("H99").Value = 2879.9092
$ws.Range("I99").Value = 2197.875
$ws.Range("J99").Value = 4698.6665
$ws.Range("K99").Value = 2197.875
$ws.Range("L99").Value = 4698.6665
$ws.Range("M99").Value = -699.875
$ws.Range("N99").Value = -7694.6665
$ws.Range("H126").Value = 2879.9092
$ws.Range("I126").Value = 2197.875
$ws.Range("J126").Value = 4698.6665
$ws.Range("K126").Value = 6593.625
$ws.Range("L126").Value = 14095.9995
$ws.Range("M126").Value = -4123.625
$ws.Range("N126").Value = -19035.9995
$ws.Range("H133").Value = 79796.42999999999
$ws.Range("J133").Value = 81431.164
$ws.Range("L133").Value = 81431.164
$ws.Range("N133").Value = -86491.164
$ws.Range("H141").Value = 422221
$ws.Range("J141").Value = 422221
$ws.Range("L141").Value = 422221
$ws.Range("N141").Value = -432581

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 34
$ws.Range("I2").Value = 38.22222
$ws.Range("K2").Value = 229.33332
$ws.Range("M2").Value = -116.33332
$ws.Range("H41").Value = 440
$ws.Range("I41").Value = 440
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1320
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -982
$ws.Range("N41").ClearContents()
$ws.Range("H122").Value = 386.88235
$ws.Range("I122").Value = 330.36365
$ws.Range("J122").Value = 490.5
$ws.Range("K122").Value = 2973.27285
$ws.Range("L122").Value = 4414.5
$ws.Range("M122").Value = -523.2728500000003
$ws.Range("N122").Value = -9314.5
$ws.Range("H123").Value = 2346.7144
$ws.Range("I123").Value = 2346.7144
$ws.Range("K123").Value = 7040.1432
$ws.Range("M123").Value = -4590.1432
$ws.Range("H129").Value = 851648.9
$ws.Range("J129").Value = 2835666.8
$ws.Range("L129").Value = 8507000.399999999
$ws.Range("N129").Value = -8517000.399999999
$ws.Range("H131").Value = 1891775.5
$ws.Range("I131").Value = 1708.3636
$ws.Range("K131").Value = 5125.0908
$ws.Range("M131").Value = -85.09079999999994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7434.625
$ws.Range("J70").Value = 6630
$ws.Range("L70").Value = 6630
$ws.Range("N70").Value = -7170
$ws.Range("H73").Value = 7434.625
$ws.Range("J73").Value = 6630
$ws.Range("L73").Value = 6630
$ws.Range("N73").Value = -8502
$ws.Range("H122").Value = 3098.45
$ws.Range("I122").Value = 3087.4443
$ws.Range("J122").Value = 3197.5
$ws.Range("K122").Value = 9262.332900000001
$ws.Range("L122").Value = 9592.5
$ws.Range("M122").Value = -6812.332900000001
$ws.Range("N122").Value = -14492.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10765.75
$ws.Range("I7").Value = 12134.77
$ws.Range("J7").Value = 4833.3335
$ws.Range("K7").Value = 12134.77
$ws.Range("L7").Value = 4833.3335
$ws.Range("M7").Value = -12022.77
$ws.Range("N7").Value = -5057.3335
$ws.Range("H126").Value = 10765.75
$ws.Range("I126").Value = 12134.77
$ws.Range("J126").Value = 4833.3335
$ws.Range("K126").Value = 36404.31
$ws.Range("L126").Value = 14500.0005
$ws.Range("M126").Value = -33934.31
$ws.Range("N126").Value = -19440.0005
$ws.Range("H132").Value = 70113.44500000001
$ws.Range("I132").Value = 78127.625
$ws.Range("K132").Value = 234382.875
$ws.Range("M132").Value = -231852.875
$ws.Range("H133").Value = 67993
$ws.Range("J133").Value = 67993
$ws.Range("L133").Value = 67993
$ws.Range("N133").Value = -73053
$ws.Range("H136").Value = 3285.5
$ws.Range("I136").Value = 2766.2
$ws.Range("K136").Value = 8298.599999999999
$ws.Range("M136").Value = -5748.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1124.0869
$ws.Range("J113").Value = 1794.8
$ws.Range("L113").Value = 5384.4
$ws.Range("N113").Value = -9724.4
$ws.Range("H126").Value = 69302.44500000001
$ws.Range("I126").Value = 73639.7
$ws.Range("J126").Value = 10749.5
$ws.Range("K126").Value = 220919.1
$ws.Range("L126").Value = 32248.5
$ws.Range("M126").Value = -218449.1
$ws.Range("N126").Value = -37188.5
$ws.Range("H132").Value = 54377.29
$ws.Range("I132").Value = 60358.652
$ws.Range("J132").Value = 5529.5
$ws.Range("K132").Value = 181075.956
$ws.Range("L132").Value = 16588.5
$ws.Range("M132").Value = -178545.956
$ws.Range("N132").Value = -21648.5
$ws.Range("H136").Value = 3993.225
$ws.Range("I136").Value = 4105.9033
$ws.Range("J136").Value = 3605.111
$ws.Range("K136").Value = 12317.7099
$ws.Range("L136").Value = 10815.333
$ws.Range("M136").Value = -9767.7099
$ws.Range("N136").Value = -15915.333
